$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$used = $ws.UsedRange
$rows = $used.Rows.Count
$cols = $used.Columns.Count

for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $cell.Value2
        if ($v -ceq "stimuli/blank.JPG") {
            $cell.Value2 = "stimuli/blank.jpg"
        }
    }
}
